# Auto-generated edit script applying cryptos.xlsx price/volume/name/link updates
# (commit: Updated cryptos list on Tue Oct  3 01:37:49 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.498.19"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "1.665.94"
$ws.Range("E3").Value = "  -3.35%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.05"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.55"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0621"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0878"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "1.898.73"
$ws.Range("E12").Value = "  -3.44%  "
$ws.Range("D13").Value = "1.657.55"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.14"
$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.28"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "248.58"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "27.491.94"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("E20").Value = "  -4.61%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.49"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.31"
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("E24").Value = "  -5.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.22"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.60"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.15"
$ws.Range("E27").Value = "  -4.75%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.112"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("E30").Value = "  +5.01%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0509"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D33").Value = "1.463.38"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  -5.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  -6.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.940"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.573"
$ws.Range("E38").Value = "  -6.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0171"
$ws.Range("E39").Value = "  -2.09%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "69.69"
$ws.Range("E40").Value = "  -2.40%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.01"
$ws.Range("E42").Value = "  -9.25%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.21"
$ws.Range("E44").Value = "  -3.12%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.807.76"
$ws.Range("E45").Value = "  -3.36%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.790"
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.69"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.54"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "0.0₆0109"
$ws.Range("E49").Value = "  -2.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "41.95"
$ws.Range("E50").Value = "  +18.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.102"
$ws.Range("E51").Value = "  -2.80%  "
